$wb = $excel.ActiveWorkbook

# --- Add the new "Spain" market sheet -------------------------------------
# The existing "Italy" sheet is the template: duplicate it and place the
# copy right after it (matching sheetId=20 / rId8 in the target workbook).
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy($null, $italy)

$spain = $wb.Worksheets.Item("Italy (2)")
$spain.Name = "Spain"

# --- Spain-specific content -------------------------------------------------
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2056"

# --- Column widths on the new sheet -----------------------------------------
# Column B narrows to fit "Spain Market" / "NGC-3103/T2056", column D narrows
# a little too.
$spain.Columns.Item(2).ColumnWidth = 29.59
$spain.Columns.Item(4).ColumnWidth = 20.1

# --- Rows 3-5 grow (text wraps because the columns are narrower now) -------
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8

# --- Selections / active sheet ----------------------------------------------
# Italy is no longer the active tab; its selection resets to the full table.
[void]$italy.Range("A1:D19").Select()

# Spain becomes the active sheet/tab, selected at C4.
[void]$spain.Activate()
[void]$spain.Range("C4").Select()
